# Rename the tool-specific worksheets to lowercase.
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("BAP").Name = "bap"
$wb.Worksheets.Item("Hopper").Name = "hopper"
$wb.Worksheets.Item("IDAPro").Name = "idapro"
$wb.Worksheets.Item("Ghidra").Name = "ghidra"
$wb.Worksheets.Item("Dyninst").Name = "dyninst"

# The chart on the Statistics sheet caches series formulas that still
# point at the old (mixed-case) sheet names; update them to match.
$ws = $wb.Worksheets.Item("Statistics")
$cos = $ws.ChartObjects()
$co1 = $cos.Item(1)
$chart1 = $co1.Chart
$series1 = $chart1.SeriesCollection()

$series1.Item(4).Formula = "=SERIES(,objdump!`$A`$2:`$A`$52,bap!`$F`$2:`$F`$52,4)"
$series1.Item(5).Formula = "=SERIES(,objdump!`$A`$2:`$A`$52,hopper!`$F`$2:`$F`$52,5)"
$series1.Item(6).Formula = "=SERIES(,objdump!`$A`$2:`$A`$52,idapro!`$F`$2:`$F`$52,6)"
$series1.Item(7).Formula = "=SERIES(,objdump!`$A`$2:`$A`$52,ghidra!`$F`$2:`$F`$52,7)"
$series1.Item(8).Formula = "=SERIES(,objdump!`$A`$2:`$A`$52,dyninst!`$F`$2:`$F`$52,8)"

# Reposition / resize the second chart on the Statistics sheet.
$co2 = $cos.Item(2)
$co2.Left = 284.343307
$co2.Top = 424.374803
$co2.Width = 764.192126
$co2.Height = 348.633071

# A new (blank) row was appended below the existing data.
$cell = $ws.Cells.Item(104, 9)
$cell.Style = "Normal"

# Update the sheet view so the new row is visible/selected.
$ws.Range("O62").Select()
